$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates continue the existing daily series)
$newRows = @(
    @{ Row = 230; A = 44304; B = 0; C = 4; D = 175.1313485113835 },
    @{ Row = 231; A = 44305; B = 0; C = 1; D = 43.78283712784589 },
    @{ Row = 232; A = 44306; B = 0; C = 1; D = 43.78283712784589 },
    @{ Row = 233; A = 44307; B = 2; C = 3; D = 131.3485113835376 }
)

# Use the last existing data row (229) as formatting template for column A (date style)
$templateCell = $ws.Range("A229")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    # Copy formatting (style) from the prior row's date cell onto the new date cell
    $templateCell.Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
